# Applies the crypto price/volume updates described by the commit diff.
# All data cells in this sheet are stored as inline strings in the source
# workbook, so every write below is a text assignment. Numeric-looking
# values (plain decimals without a thousands separator) are written with a
# leading apostrophe to force Excel to keep them as text instead of auto-
# converting them to numbers; the style is then reset to "Normal" so the
# cell keeps the workbook's original (unformatted) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.589.21"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.891.95"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'238.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.4886"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "'0.2929"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.52%  "
$ws.Range("D9").Value = "'0.06692"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.53%  "
$ws.Range("D10").Value = "1.898.43"
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'17.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "'0.07338"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").Value = "'5.150"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.02%  "
$ws.Range("D14").Value = "'88.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("D15").Value = "'0.6665"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("D16").Value = "30.541.03"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.61%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.000007834"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("E19").Value = "  -0.20%  "
$ws.Range("D20").Value = "2.139.61"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").Value = "'1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.282"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +11.57%  "
$ws.Range("D23").Value = "'189.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "'6.162"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.08%  "
$ws.Range("D25").Value = "'9.461"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").Value = "'161.77"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").Value = "'18.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'1.928"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.31%  "
$ws.Range("D29").Value = "'1.476"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.97%  "
$ws.Range("D30").Value = "'4.364"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.59%  "
$ws.Range("D31").Value = "'0.09141"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.26%  "
$ws.Range("D32").Value = "'4.106"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.55%  "
$ws.Range("D33").Value = "'0.05212"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").Value = "'0.7366"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").Value = "'1.098"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "'2.715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("D37").Value = "'0.01829"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "'2.684"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").Value = "'0.9208"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "'2.050"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "'0.4397"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Value = "'5.925"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.68%  "
$ws.Range("D43").Value = "'106.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").Value = "'0.9933"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'68.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +20.93%  "
$ws.Range("D46").Value = "'0.1382"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.21%  "
$ws.Range("D47").Value = "'7.566"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D48").Value = "'8.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.19%  "
$ws.Range("D49").Value = "'34.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("D50").Value = "'0.05827"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'1.422"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.12%  "
